$d = $word.ActiveDocument

# 1. Replace the blank underscore placeholder with "(e)" in the transfer-tax
#    exemption paragraph (appears in both the drawing and VML fallback copies
#    of the text box, and Word keeps both in sync on save).
$d.Content.Find.Execute("Paragraph ___, Illinois Real Estate Transfer Tax Law.", $true, $false, $false, $false, $false, $true, 1, $false, "Paragraph (e), Illinois Real Estate Transfer Tax Law.", 2)

# 2. Remove the now-stale _GoBack bookmark left over from the last edit
#    position, and make sure the "in my presence on" text reads as a single
#    run-worth of text.
$d.Content.Find.Execute("executed this Illinois Transfer on Death Instrument in my presence on", $true, $false, $false, $false, $false, $true, 1, $false, "executed this Illinois Transfer on Death Instrument in my presence on", 2)
